# Update the "Förändrad" (Changed) date column (C) for rows 2-16
# from serial date 45174 (2023-09-05) to 45175 (2023-09-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value()
    if ($current -ne $null -and $current.ToOADate() -eq 45174) {
        $cell.Value = 45175
    }
}
